$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "'90"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'1"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'1"
$ws.Range("G3").Style = "Normal"
$ws.Range("J4").Value = "'6"
$ws.Range("J4").Style = "Normal"
$ws.Range("E7").Value = "'900"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'10"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'10"
$ws.Range("G7").Style = "Normal"
$ws.Range("L7").Value = "'5"
$ws.Range("L7").Style = "Normal"
$ws.Range("E8").Value = "'743"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'11"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'7"
$ws.Range("G8").Style = "Normal"
$ws.Range("J9").Value = "'7"
$ws.Range("J9").Style = "Normal"
$ws.Range("E10").Value = "'417"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'7"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'5"
$ws.Range("G10").Style = "Normal"
$ws.Range("L10").Value = "'2"
$ws.Range("L10").Style = "Normal"
$ws.Range("E11").Value = "'185"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'3"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'2"
$ws.Range("G11").Style = "Normal"
$ws.Range("I11").Value = "'1"
$ws.Range("I11").Style = "Normal"
$ws.Range("E12").Value = "'990"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'11"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "'11"
$ws.Range("G12").Style = "Normal"
$ws.Range("L12").Value = "'2"
$ws.Range("L12").Style = "Normal"
$ws.Range("E14").Value = "'14"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = "'2"
$ws.Range("F14").Style = "Normal"
$ws.Range("H14").Value = "'2"
$ws.Range("H14").Style = "Normal"
$ws.Range("J14").Value = "'2"
$ws.Range("J14").Style = "Normal"
$ws.Range("J15").Value = "'7"
$ws.Range("J15").Style = "Normal"
$ws.Range("E16").Value = "'556"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = "'11"
$ws.Range("F16").Style = "Normal"
$ws.Range("H16").Value = "'4"
$ws.Range("H16").Style = "Normal"
$ws.Range("J16").Value = "'4"
$ws.Range("J16").Style = "Normal"
$ws.Range("E17").Value = "'352"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = "'11"
$ws.Range("F17").Style = "Normal"
$ws.Range("H17").Value = "'9"
$ws.Range("H17").Style = "Normal"
$ws.Range("J17").Value = "'9"
$ws.Range("J17").Style = "Normal"
$ws.Range("E18").Value = "'607"
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").Value = "'10"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = "'8"
$ws.Range("G18").Style = "Normal"
$ws.Range("I18").Value = "'8"
$ws.Range("I18").Style = "Normal"
$ws.Range("E19").Value = "'692"
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").Value = "'10"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = "'8"
$ws.Range("G19").Style = "Normal"
$ws.Range("J20").Value = "'9"
$ws.Range("J20").Style = "Normal"
$ws.Range("E23").Value = "'186"
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").Value = "'7"
$ws.Range("F23").Style = "Normal"
$ws.Range("H23").Value = "'5"
$ws.Range("H23").Style = "Normal"
$ws.Range("J23").Value = "'9"
$ws.Range("J23").Style = "Normal"
$ws.Range("E25").Value = "'915"
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").Value = "'11"
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").Value = "'11"
$ws.Range("G25").Style = "Normal"
$ws.Range("I25").Value = "'6"
$ws.Range("I25").Style = "Normal"
$ws.Range("E26").Value = "'866"
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").Value = "'11"
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").Value = "'10"
$ws.Range("G26").Style = "Normal"
$ws.Range("I26").Value = "'5"
$ws.Range("I26").Style = "Normal"
$ws.Range("E27").Value = "'696"
$ws.Range("E27").Style = "Normal"
$ws.Range("F27").Value = "'11"
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").Value = "'9"
$ws.Range("G27").Style = "Normal"
$ws.Range("I27").Value = "'7"
$ws.Range("I27").Style = "Normal"
$ws.Range("E29").Value = "'291"
$ws.Range("E29").Style = "Normal"
$ws.Range("F29").Value = "'9"
$ws.Range("F29").Style = "Normal"
$ws.Range("H29").Value = "'7"
$ws.Range("H29").Style = "Normal"
$ws.Range("J29").Value = "'8"
$ws.Range("J29").Style = "Normal"
